$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Years first
$ws.Range("A15").Value = 1991
$ws.Range("A16").Value = 1991

# Japanese titles (column B) for both new rows
$ws.Range("B15").Value = "ビヨンド・ローズ・トゥ・ロードでわかる実践RPG入門"
$ws.Range("B16").Value = "精霊の大地―ビヨンド・ローズ・トゥ・ロードでわかる実践RPG入門 2"

# Publisher (column D) for both new rows
$ws.Range("D15").Value = "BNN"
$ws.Range("D16").Value = "BNN"

# Image (column E) for both new rows
$ws.Range("E15").Value = "beyond_roads_to_lord_primer1.jpg"
$ws.Range("E16").Value = "beyond_roads_to_lord_primer2.jpg"

# English titles (column C) for both new rows
$ws.Range("C15").Value = "A practical introduction to RPGs by means of Beyond Roads to Lord"
$ws.Range("C16").Value = "Earth Spirits: A practical introduction to RPGs by means of Beyond Roads to Lord Vol 2"

# Product type (column F) for both new rows
$ws.Range("F15").Value = "supplement"
$ws.Range("F16").Value = "supplement"

# Column C width change (48.83203125 -> 61.33203125)
$ws.Range("C1").ColumnWidth = 60.41666666666667

# Update selection to C17
$ws.Range("C17").Select()
